$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 186
$ws.Range("F6").Value = 167
$ws.Range("F7").Value = 159
$ws.Range("F8").Value = 783
$ws.Range("F9").Value = 4164
$ws.Range("F11").Value = 56
$ws.Range("F12").Value = 171
$ws.Range("F14").Value = 5985
$ws.Range("F15").Value = 463
$ws.Range("F16").Value = 2307
$ws.Range("F19").Value = 458
$ws.Range("F20").Value = 9005
$ws.Range("F22").Value = 2331
$ws.Range("G22").Value = 44.1
$ws.Range("F23").Value = 195
$ws.Range("F24").Value = 2292
$ws.Range("F25").Value = 2408
$ws.Range("F26").Value = 1378
$ws.Range("F27").Value = 229
$ws.Range("F28").Value = 1940
$ws.Range("F30").Value = 55
$ws.Range("F34").Value = 39
$ws.Range("F35").Value = 40
$ws.Range("F36").Value = 32
$ws.Range("F37").Value = 1221
$ws.Range("F39").Value = 70
$ws.Range("F40").Value = 91
$ws.Range("F41").Value = 229
$ws.Range("F42").Value = 1508
$ws.Range("F43").Value = 2459
$ws.Range("F45").Value = 908
$ws.Range("F46").Value = 288
$ws.Range("F47").Value = 1246
$ws.Range("F48").Value = 10

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 932
$ws.Range("F12").Value = 146
$ws.Range("F21").Value = 27
$ws.Range("F22").Value = 37
$ws.Range("F23").Value = 37

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 686
$ws.Range("F3").Value = 883
$ws.Range("F4").Value = 97

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 686
$ws.Range("F4").Value = 883
$ws.Range("F5").Value = 97
$ws.Range("F8").Value = 167
$ws.Range("F11").Value = 159
$ws.Range("F12").Value = 783
$ws.Range("F13").Value = 4164
$ws.Range("F14").Value = 4164
$ws.Range("F15").Value = 56
$ws.Range("F16").Value = 171
$ws.Range("F19").Value = 5985
$ws.Range("F20").Value = 463
$ws.Range("F21").Value = 2307
$ws.Range("F23").Value = 458
$ws.Range("F24").Value = 9005
$ws.Range("F25").Value = 146
$ws.Range("F27").Value = 2331
$ws.Range("G27").Value = 44.1
$ws.Range("F28").Value = 2292
$ws.Range("F29").Value = 2408
$ws.Range("F30").Value = 1378
$ws.Range("F31").Value = 229
$ws.Range("F32").Value = 1940
$ws.Range("F34").Value = 55
$ws.Range("F37").Value = 40
$ws.Range("F38").Value = 32
$ws.Range("F39").Value = 1221
$ws.Range("F41").Value = 91
$ws.Range("F42").Value = 229
$ws.Range("F43").Value = 1508
$ws.Range("F44").Value = 2459
$ws.Range("F45").Value = 908
$ws.Range("F46").Value = 288
$ws.Range("F49").Value = 27
$ws.Range("F50").Value = 1246
$ws.Range("F51").Value = 37
